# New crime data collected — weekly refresh of the 72nd Precinct CompStat
# report: bump the volume/number + reporting week header, then overwrite
# the weekly/28-day/YTD/2-year crime-count grid (C14:N33) with the newly
# collected figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text: "Volume 32   Number  26" -> "...27", and the reporting
# week date range moves forward one week.
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 32   Number  27"
$ws.Range("C9").Value = "Report Covering the Week  6/30/2025  Through  7/6/2025"

# ---------------------------------------------------------------------
# Helpers.
#
# The data grid stores "no data" placeholders as literal text ("0" or
# "***.*") sharing the same look as numbers. Excel's COM layer infers a
# numeric type from a plain numeric-looking string, so any cell that
# must flip from a number to one of these text placeholders needs its
# number format forced to Text ("@") first; the General numeric cells
# conversely need their numeric format restored when a placeholder is
# replaced by a real number.
# ---------------------------------------------------------------------
function Set-NumCell($addr, $value) {
    $ws.Range($addr).Value = $value
}

function Set-TextCell($addr, $text) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
}

$FMT_INT = "#,##0"
$FMT_PCT1 = "#,##0.0;""-""#,##0.0"

function Set-IntCell($addr, $value) {
    $ws.Range($addr).NumberFormat = $FMT_INT
    $ws.Range($addr).Value = $value
}

function Set-PctCell($addr, $value) {
    $ws.Range($addr).NumberFormat = $FMT_PCT1
    $ws.Range($addr).Value = $value
}

# ---------------------------------------------------------------------
# Row 14 - Murder: WTD 2025 % chg column flips from a real number back
# to the "no activity" text placeholders.
# ---------------------------------------------------------------------
Set-TextCell "G14" "0"
Set-TextCell "H14" "***.*"

# ---------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------
Set-TextCell "C15" "0"
Set-NumCell  "D15" 2
Set-NumCell  "E15" -100
Set-NumCell  "F15" 1
Set-NumCell  "G15" 4
Set-NumCell  "H15" -75
Set-NumCell  "J15" 14
Set-NumCell  "K15" -21.428571428571
Set-NumCell  "N15" -38.888888888888

# ---------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------
Set-NumCell "D16" 4
Set-NumCell "E16" -50
Set-NumCell "F16" 6
Set-NumCell "G16" 22
Set-NumCell "H16" -72.727272727272
Set-NumCell "I16" 48
Set-NumCell "J16" 101
Set-NumCell "K16" -52.475247524752
Set-NumCell "L16" -32.394366197183
Set-NumCell "M16" -49.473684210526
Set-NumCell "N16" -91.504424778761

# ---------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------
Set-NumCell "C17" 10
Set-NumCell "D17" 5
Set-NumCell "E17" 100
Set-NumCell "F17" 35
Set-NumCell "G17" 26
Set-NumCell "H17" 34.615384615384
Set-NumCell "I17" 146
Set-NumCell "J17" 166
Set-NumCell "K17" -12.048192771084
Set-NumCell "L17" -1.351351351351
Set-NumCell "M17" 56.989247311828
Set-NumCell "N17" -48.771929824561

# ---------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------
Set-NumCell "D18" 2
Set-NumCell "E18" 0
Set-NumCell "F18" 6
Set-NumCell "G18" 10
Set-NumCell "H18" -40
Set-NumCell "I18" 41
Set-NumCell "J18" 87
Set-NumCell "K18" -52.873563218390
Set-NumCell "L18" -36.923076923076
Set-NumCell "M18" -69.172932330827
Set-NumCell "N18" -92.979452054794

# ---------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------
Set-NumCell "C19" 3
Set-NumCell "D19" 12
Set-NumCell "E19" -75
Set-NumCell "F19" 29
Set-NumCell "G19" 41
Set-NumCell "H19" -29.268292682926
Set-NumCell "I19" 170
Set-NumCell "J19" 261
Set-NumCell "K19" -34.865900383141
Set-NumCell "L19" -39.068100358422
Set-NumCell "M19" -1.734104046242
Set-NumCell "N19" -24.107142857142

# ---------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------
Set-NumCell "C20" 3
Set-NumCell "D20" 2
Set-NumCell "E20" 50
Set-NumCell "F20" 11
Set-NumCell "G20" 14
Set-NumCell "H20" -21.428571428571
Set-NumCell "I20" 39
Set-NumCell "J20" 84
Set-NumCell "K20" -53.571428571428
Set-NumCell "L20" -53.571428571428
Set-NumCell "M20" -35
Set-NumCell "N20" -92.3828125

# ---------------------------------------------------------------------
# Row 21 - TOTAL (bold row, keeps its own style family already)
# ---------------------------------------------------------------------
Set-NumCell "C21" 20
Set-NumCell "D21" 27
Set-NumCell "E21" -25.925925925925
Set-NumCell "F21" 88
Set-NumCell "G21" 117
Set-NumCell "H21" -24.786324786324
Set-NumCell "I21" 455
Set-NumCell "J21" 714
Set-NumCell "K21" -36.274509803921
Set-NumCell "L21" -30.851063829787
Set-NumCell "M21" -20.035149384885
Set-NumCell "N21" -79.299363057324

# ---------------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------------
Set-TextCell "C22" "0"
Set-NumCell  "F22" 5
Set-NumCell  "H22" 150
Set-NumCell  "L22" 18.75

# Row 23 - Housing: unchanged in this update.

# ---------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------
Set-NumCell "C24" 13
Set-NumCell "D24" 13
Set-NumCell "E24" 0
Set-NumCell "F24" 60
Set-NumCell "G24" 97
Set-NumCell "H24" -38.144329896907
Set-NumCell "I24" 399
Set-NumCell "J24" 510
Set-NumCell "K24" -21.764705882352
Set-NumCell "L24" -33.388981636060
Set-NumCell "M24" 7.547169811320

# ---------------------------------------------------------------------
# Row 25 - Retail Theft
# ---------------------------------------------------------------------
Set-NumCell "C25" 9
Set-NumCell "D25" 8
Set-NumCell "E25" 12.5
Set-NumCell "F25" 18
Set-NumCell "G25" 33
Set-NumCell "H25" -45.454545454545
Set-NumCell "I25" 125
Set-NumCell "J25" 229
Set-NumCell "K25" -45.414847161572
Set-NumCell "L25" -61.300309597523

# ---------------------------------------------------------------------
# Row 26 - Misd. Assault
# ---------------------------------------------------------------------
Set-NumCell "C26" 15
Set-NumCell "D26" 13
Set-NumCell "E26" 15.384615384615
Set-NumCell "F26" 48
Set-NumCell "G26" 35
Set-NumCell "H26" 37.142857142857
Set-NumCell "I26" 260
Set-NumCell "J26" 248
Set-NumCell "K26" 4.838709677419
Set-NumCell "L26" 1.960784313725
Set-NumCell "M26" -22.155688622754

# ---------------------------------------------------------------------
# Row 27 - UCR Rape*
# ---------------------------------------------------------------------
Set-TextCell "C27" "0"
Set-NumCell  "D27" 2
Set-NumCell  "E27" -100
Set-NumCell  "F27" 1
Set-NumCell  "G27" 4
Set-NumCell  "H27" -75
Set-NumCell  "J27" 19
Set-NumCell  "K27" -15.789473684210

# ---------------------------------------------------------------------
# Row 28 - Other Sex Crimes
# ---------------------------------------------------------------------
Set-NumCell "C28" 3
Set-IntCell "D28" 1
Set-PctCell "E28" 200
Set-NumCell "F28" 7
Set-NumCell "G28" 3
Set-NumCell "H28" 133.333333333333
Set-NumCell "I28" 32
Set-NumCell "J28" 29
Set-NumCell "K28" 10.344827586206
Set-NumCell "L28" -30.434782608695

# Rows 29/30 - Shooting Vic. / Shooting Inc.: unchanged in this update.

# ---------------------------------------------------------------------
# Row 31 - Hate Crimes
# ---------------------------------------------------------------------
Set-TextCell "D31" "0"
Set-TextCell "E31" "***.*"
Set-TextCell "F31" "0"
Set-NumCell  "H31" -100

# ---------------------------------------------------------------------
# Row 33 - Traffic Fatalities
# ---------------------------------------------------------------------
Set-TextCell "G33" "0"
Set-TextCell "H33" "***.*"
